$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue "D2" "57.112.81"
Set-TextValue "E2" "  -1.09%  "
Set-TextValue "D3" "2.986.57"
Set-TextValue "E3" "  -2.06%  "
Set-TextValue "E4" "  -0.08%  "
Set-TextValue "D5" "501.30"
Set-TextValue "E5" "  -4.32%  "
Set-TextValue "D6" "137.90"
Set-TextValue "E6" "  -2.99%  "
Set-TextValue "E7" "  +0.04%  "
Set-TextValue "E8" "  -3.93%  "
Set-TextValue "E9" "  -4.27%  "
Set-TextValue "E10" "  -4.28%  "
Set-TextValue "E11" "  -3.56%  "
Set-TextValue "D12" "3.500.12"
Set-TextValue "E12" "  -2.15%  "
Set-TextValue "E13" "  -2.36%  "
Set-TextValue "D14" "26.08"
Set-TextValue "E14" "  -3.32%  "
Set-TextValue "E15" "  -5.73%  "
Set-TextValue "D16" "57.157.58"
Set-TextValue "E16" "  -0.99%  "
Set-TextValue "D17" "6.07"
Set-TextValue "E17" "  -3.45%  "
Set-TextValue "D18" "2.985.87"
Set-TextValue "E18" "  -2.26%  "
Set-TextValue "E19" "  -3.33%  "
Set-TextValue "E20" "  -3.48%  "
Set-TextValue "D21" "321.42"
Set-TextValue "E21" "  -4.97%  "
Set-TextValue "E22" "  -0.06%  "
Set-TextValue "D23" "5.71"
Set-TextValue "E23" "  +0.52%  "
Set-TextValue "D24" "0.492"
Set-TextValue "E24" "  -1.81%  "
Set-TextValue "D25" "63.12"
Set-TextValue "E25" "  -2.77%  "
Set-TextValue "D26" "1.01"
Set-TextValue "E26" "  +0.43%  "
Set-TextValue "E27" "  -5.36%  "
Set-TextValue "E28" "  -8.69%  "
Set-TextValue "D29" "6.65"
Set-TextValue "E29" "  -3.70%  "
Set-TextValue "E31" "  -4.13%  "
Set-TextValue "D32" "1.16"
Set-TextValue "E32" "  -5.11%  "
Set-TextValue "D33" "20.23"
Set-TextValue "E33" "  -4.25%  "
Set-TextValue "D34" "155.31"
Set-TextValue "E34" "  -0.87%  "
Set-TextValue "E35" "  -3.49%  "
Set-TextValue "D36" "5.79"
Set-TextValue "E37" "  -6.31%  "
Set-TextValue "D38" "24.46"
Set-TextValue "E38" "  -6.32%  "
Set-TextValue "D39" "0.0665"
Set-TextValue "E39" "  -5.74%  "
Set-TextValue "D40" "37.85"
Set-TextValue "E40" "  +0.30%  "
Set-TextValue "D41" "3.018.07"
Set-TextValue "E41" "  -2.21%  "
Set-TextValue "D42" "1.00"
Set-TextValue "E42" "  -0.01%  "
Set-TextValue "E43" "  -3.53%  "
Set-TextValue "E44" "  -2.79%  "
Set-TextValue "D45" "2.193.36"
Set-TextValue "E45" "  -5.86%  "
Set-TextValue "E46" "  -6.02%  "
Set-TextValue "E47" "  -1.77%  "
Set-TextValue "D48" "0.935"
Set-TextValue "E48" "  -9.61%  "
Set-TextValue "D49" "0.0234"
Set-TextValue "E49" "  -4.99%  "
Set-TextValue "D50" "19.29"
Set-TextValue "E50" "  -4.32%  "
Set-TextValue "E51" "  -11.08%  "
